$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "id" column (D) with a sequential row id for every data row.
$ws.Range("D1").Value = "id"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 164 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $r - 1
}

# Keep the last-used selection in the same place the author left it.
$ws.Range("E163").Select()
